$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partnership_Emails")

$row = 24
$srcRow = $row - 1

$ws.Range("A$srcRow").Copy()
$ws.Range("A$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 45857.00277777778
$ws.Cells.Item($row, 2).Value = "Moris Mwai"
$ws.Cells.Item($row, 3).Value = "Tech-Neo GmbH"
$ws.Cells.Item($row, 4).Value = "Am main City, Germany"
$ws.Cells.Item($row, 5).Value = "DE1567890"
$ws.Cells.Item($row, 6).Value = "morismwai1@gmail.com"
$ws.Cells.Item($row, 7).Value = "Partnership Offer"
